# Update visitor-count ("想去人数") figures (column F) across the three
# data sheets that share the same underlying rows: 展览, 演出 and 全部类型.
# The sheet 本地生活 has no data rows and is left untouched.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 63
$ws1.Range("F3").Value  = 583
$ws1.Range("F4").Value  = 43
$ws1.Range("F6").Value  = 10
$ws1.Range("F7").Value  = 14597
$ws1.Range("F9").Value  = 666
$ws1.Range("F10").Value = 15151
$ws1.Range("F12").Value = 8578
$ws1.Range("F13").Value = 311
$ws1.Range("F15").Value = 58
$ws1.Range("F16").Value = 174
$ws1.Range("F18").Value = 179
$ws1.Range("F19").Value = 9
$ws1.Range("F20").Value = 5
$ws1.Range("F21").Value = 19
$ws1.Range("F24").Value = 1068
$ws1.Range("F25").Value = 4
$ws1.Range("F27").Value = 49
$ws1.Range("F28").Value = 27
$ws1.Range("F30").Value = 413
$ws1.Range("F31").Value = 19
$ws1.Range("F34").Value = 262
$ws1.Range("F35").Value = 414
$ws1.Range("F37").Value = 5294
$ws1.Range("F38").Value = 5225

# --- Sheet "演出" ---------------------------------------------------------
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 1002

# --- Sheet "全部类型" ------------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 63
$ws4.Range("F3").Value  = 583
$ws4.Range("F4").Value  = 43
$ws4.Range("F6").Value  = 10
$ws4.Range("F7").Value  = 14597
$ws4.Range("F9").Value  = 666
$ws4.Range("F10").Value = 15151
$ws4.Range("F12").Value = 8578
$ws4.Range("F13").Value = 311
$ws4.Range("F15").Value = 1002
$ws4.Range("F16").Value = 58
$ws4.Range("F17").Value = 174
$ws4.Range("F19").Value = 179
$ws4.Range("F20").Value = 9
$ws4.Range("F21").Value = 5
$ws4.Range("F22").Value = 19
$ws4.Range("F25").Value = 1068
$ws4.Range("F26").Value = 4
$ws4.Range("F28").Value = 49
$ws4.Range("F29").Value = 27
$ws4.Range("F33").Value = 413
$ws4.Range("F34").Value = 19
$ws4.Range("F37").Value = 262
$ws4.Range("F38").Value = 414
$ws4.Range("F40").Value = 5294
$ws4.Range("F41").Value = 5225
